$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.162423666666667
$ws.Range("H2").Value = 3.487271
$ws.Range("I2").Value = 0.6447270069705344
$ws.Range("J2").Value = 0.6447270069705344
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.672264666666667
$ws.Range("N2").Value = 8.016794000000001
$ws.Range("O2").Value = 0.06772620019093417
$ws.Range("P2").Value = 0.06772620019093417
$ws.Range("Q2").Value = 3.106303692130445
$ws.Range("R2").Value = 27.95673322917401
$ws.Range("S2").Value = 0.04366491034258822
$ws.Range("T2").Value = 0.04366491034258822
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.162423666666667
$ws.Range("H3").Value = 3.487271
$ws.Range("I3").Value = 0.6447270069705344
$ws.Range("J3").Value = 0.6447270069705344
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.41886
$ws.Range("N3").Value = 82.25658
$ws.Range("O3").Value = 0.6949069171668364
$ws.Range("P3").Value = 0.6949069171668364
$ws.Range("Q3").Value = 31.87233177702
$ws.Range("R3").Value = 286.85098599318
$ws.Range("S3").Value = 0.4480252568280955
$ws.Range("T3").Value = 0.4480252568280955
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.162423666666667
$ws.Range("H4").Value = 3.487271
$ws.Range("I4").Value = 0.6447270069705344
$ws.Range("J4").Value = 0.6447270069705344
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 9.365757
$ws.Range("N4").Value = 28.097271
$ws.Range("O4").Value = 0.2373668826422294
$ws.Range("P4").Value = 0.2373668826422294
$ws.Range("Q4").Value = 10.886977593049
$ws.Range("R4").Value = 97.982798337441
$ws.Range("S4").Value = 0.1530368397998507
$ws.Range("T4").Value = 0.1530368397998507
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.4242653333333333
$ws.Range("H5").Value = 1.272796
$ws.Range("I5").Value = 0.2353146502133239
$ws.Range("J5").Value = 0.2353146502133239
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.672264666666667
$ws.Range("N5").Value = 8.016794000000001
$ws.Range("O5").Value = 0.06772620019093417
$ws.Range("P5").Value = 0.06772620019093417
$ws.Range("Q5").Value = 1.133749259558222
$ws.Range("R5").Value = 10.203743336024
$ws.Range("S5").Value = 0.01593696710820722
$ws.Range("T5").Value = 0.01593696710820723
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.4242653333333333
$ws.Range("H6").Value = 1.272796
$ws.Range("I6").Value = 0.2353146502133239
$ws.Range("J6").Value = 0.2353146502133239
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 27.41886
$ws.Range("N6").Value = 82.25658
$ws.Range("O6").Value = 0.6949069171668364
$ws.Range("P6").Value = 0.6949069171668364
$ws.Range("Q6").Value = 11.63287177752
$ws.Range("R6").Value = 104.69584599768
$ws.Range("S6").Value = 0.1635217781439333
$ws.Range("T6").Value = 0.1635217781439334
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.4242653333333333
$ws.Range("H7").Value = 1.272796
$ws.Range("I7").Value = 0.2353146502133239
$ws.Range("J7").Value = 0.2353146502133239
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 9.365757
$ws.Range("N7").Value = 28.097271
$ws.Range("O7").Value = 0.2373668826422294
$ws.Range("P7").Value = 0.2373668826422294
$ws.Range("Q7").Value = 3.973566015524
$ws.Range("R7").Value = 35.762094139716
$ws.Range("S7").Value = 0.05585590496118332
$ws.Range("T7").Value = 0.05585590496118333
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.2162813333333333
$ws.Range("H8").Value = 0.648844
$ws.Range("I8").Value = 0.1199583428161417
$ws.Range("J8").Value = 0.1199583428161417
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.672264666666667
$ws.Range("N8").Value = 8.016794000000001
$ws.Range("O8").Value = 0.06772620019093417
$ws.Range("P8").Value = 0.06772620019093417
$ws.Range("Q8").Value = 0.5779609651262223
$ws.Range("R8").Value = 5.201648686136
$ws.Range("S8").Value = 0.008124322740138724
$ws.Range("T8").Value = 0.008124322740138724
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.2162813333333333
$ws.Range("H9").Value = 0.648844
$ws.Range("I9").Value = 0.1199583428161417
$ws.Range("J9").Value = 0.1199583428161417
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 27.41886
$ws.Range("N9").Value = 82.25658
$ws.Range("O9").Value = 0.6949069171668364
$ws.Range("P9").Value = 0.6949069171668364
$ws.Range("Q9").Value = 5.930187599279999
$ws.Range("R9").Value = 53.37168839352
$ws.Range("S9").Value = 0.08335988219480757
$ws.Range("T9").Value = 0.08335988219480757
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.2162813333333333
$ws.Range("H10").Value = 0.648844
$ws.Range("I10").Value = 0.1199583428161417
$ws.Range("J10").Value = 0.1199583428161417
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 9.365757
$ws.Range("N10").Value = 28.097271
$ws.Range("O10").Value = 0.2373668826422294
$ws.Range("P10").Value = 0.2373668826422294
$ws.Range("Q10").Value = 2.025638411636
$ws.Range("R10").Value = 18.230745704724
$ws.Range("S10").Value = 0.02847413788119544
$ws.Range("T10").Value = 0.02847413788119544
